$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column I (9th column: NeedTar), shifting
# NeedTar, DefaultHitTime, ShowName one column to the right.
$ws.Columns.Item(9).Insert()

# New header for the inserted column
$ws.Range("I1").Value = "AutoAtkDis"

# New column values for rows 2-9 (skills auto-attack within distance)
for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 9).Value = 1
}

# Match the selection left behind by the edit
$ws.Range("I2:I9").Select() | Out-Null
